$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 42
$ws.Range("H42").Value = 72.5
$ws.Range("J42").Value = 69.5
$ws.Range("L42").Value = 208.5
$ws.Range("N42").Value = -668.5

# Row 64
$ws.Range("H64").Value = 148999.47
$ws.Range("J64").Value = 8750
$ws.Range("L64").Value = 8750
$ws.Range("N64").Value = -9246

# Row 67
$ws.Range("H67").Value = 148999.47
$ws.Range("J67").Value = 8750
$ws.Range("L67").Value = 8750
$ws.Range("N67").Value = -10466

# Row 97
$ws.Range("H97").Value = 3226.75
$ws.Range("J97").Value = 2973.4285
$ws.Range("L97").Value = 8920.2855
$ws.Range("N97").Value = -9912.2855

# Row 107
$ws.Range("H107").Value = 5033.909
$ws.Range("I107").Value = 4910.8887
$ws.Range("K107").Value = 4910.8887
$ws.Range("M107").Value = -2990.8887

# Row 115
$ws.Range("H115").Value = 1574.7273
$ws.Range("I115").Value = 1216.5
$ws.Range("K115").Value = 3649.5
$ws.Range("M115").Value = -2082.5

# Row 132
$ws.Range("H132").Value = 2860790
$ws.Range("I132").Value = 3595.0667
$ws.Range("K132").Value = 10785.2001
$ws.Range("M132").Value = -8255.2001

# Row 137
$ws.Range("H137").Value = 775647.7
$ws.Range("I137").Value = 1277369.1
$ws.Range("K137").Value = 3832107.3
$ws.Range("M137").Value = -3829557.3

# Row 138
$ws.Range("H138").Value = 4047.92
$ws.Range("I138").Value = 1591.8572
$ws.Range("J138").Value = 5003.0557
$ws.Range("K138").Value = 4775.571599999999
$ws.Range("L138").Value = 15009.1671
$ws.Range("M138").Value = 364.4284000000007
$ws.Range("N138").Value = -25289.1671

$ws = $wb.Worksheets.Item("ARM")
# Row 43
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = $null

# Row 45
$ws.Range("H45").Value = 7370.7334
$ws.Range("I45").Value = 7229.1113
$ws.Range("K45").Value = 7229.1113
$ws.Range("M45").Value = -6852.1113

# Row 132
$ws.Range("H132").Value = 2986.5
$ws.Range("I132").Value = 1981.1
$ws.Range("J132").Value = 5500
$ws.Range("K132").Value = 5943.299999999999
$ws.Range("L132").Value = 16500
$ws.Range("M132").Value = -3413.299999999999
$ws.Range("N132").Value = -21560

# Row 139
$ws.Range("H139").Value = 221052
$ws.Range("J139").Value = 221052
$ws.Range("L139").Value = 221052
$ws.Range("N139").Value = -231332

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -127
$ws.Range("N22").Value = $null

# Row 81
$ws.Range("H81").Value = 130156
$ws.Range("J81").Value = 130156
$ws.Range("L81").Value = 130156
$ws.Range("N81").Value = -132278

# Row 82
$ws.Range("H82").Value = 29442.084
$ws.Range("I82").Value = 7622
$ws.Range("J82").Value = 59990.2
$ws.Range("K82").Value = 7622
$ws.Range("L82").Value = 59990.2
$ws.Range("M82").Value = -7239
$ws.Range("N82").Value = -60756.2

# Row 84
$ws.Range("H84").Value = 130156
$ws.Range("J84").Value = 130156
$ws.Range("L84").Value = 390468
$ws.Range("N84").Value = -401076

# Row 85
$ws.Range("H85").Value = 29442.084
$ws.Range("I85").Value = 7622
$ws.Range("J85").Value = 59990.2
$ws.Range("K85").Value = 7622
$ws.Range("L85").Value = 59990.2
$ws.Range("M85").Value = -6296
$ws.Range("N85").Value = -62642.2

# Row 94
$ws.Range("H94").Value = 3159.7942
$ws.Range("J94").Value = 7168
$ws.Range("L94").Value = 7168
$ws.Range("N94").Value = -8070

# Row 107
$ws.Range("H107").Value = 2082
$ws.Range("I107").Value = 2082
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2082
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -162
$ws.Range("N107").Value = $null

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 737.75
$ws.Range("J22").Value = 1495
$ws.Range("L22").Value = 1495
$ws.Range("N22").Value = -2195

# Row 87
$ws.Range("H87").Value = 45490.5
$ws.Range("J87").Value = 45490.5
$ws.Range("L87").Value = 45490.5
$ws.Range("N87").Value = -47862.5

# Row 90
$ws.Range("H90").Value = 45490.5
$ws.Range("J90").Value = 45490.5
$ws.Range("L90").Value = 136471.5
$ws.Range("N90").Value = -148327.5

# Row 132
$ws.Range("H132").Value = 8551.714
$ws.Range("I132").Value = 13487.75
$ws.Range("J132").Value = 1970.3334
$ws.Range("K132").Value = 40463.25
$ws.Range("L132").Value = 5911.0002
$ws.Range("M132").Value = -37933.25
$ws.Range("N132").Value = -10971.0002

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 1397.375
$ws.Range("I131").Value = 853.5333000000001
$ws.Range("J131").Value = 2303.7778
$ws.Range("K131").Value = 2560.5999
$ws.Range("L131").Value = 6911.3334
$ws.Range("M131").Value = 2479.4001
$ws.Range("N131").Value = -16991.3334

$ws = $wb.Worksheets.Item("GSM")
# Row 35
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").Value = $null

# Row 132
$ws.Range("H132").Value = 3021.0908
$ws.Range("J132").Value = 2372.2
$ws.Range("L132").Value = 7116.599999999999
$ws.Range("N132").Value = -12176.6

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 708.54285
$ws.Range("I22").Value = 572.8570999999999
$ws.Range("J22").Value = 912.0714
$ws.Range("K22").Value = 572.8570999999999
$ws.Range("L22").Value = 912.0714
$ws.Range("M22").Value = -277.8570999999999
$ws.Range("N22").Value = -1502.0714

# Row 27
$ws.Range("H27").Value = 708.54285
$ws.Range("I27").Value = 572.8570999999999
$ws.Range("J27").Value = 912.0714
$ws.Range("K27").Value = 572.8570999999999
$ws.Range("L27").Value = 912.0714
$ws.Range("M27").Value = -465.8570999999999
$ws.Range("N27").Value = -1126.0714

# Row 46
$ws.Range("H46").Value = 1343.8422
$ws.Range("J46").Value = 1848.3
$ws.Range("L46").Value = 1848.3
$ws.Range("N46").Value = -2224.3

# Row 136
$ws.Range("H136").Value = 4864.6924
$ws.Range("I136").Value = 3560.7273
$ws.Range("J136").Value = 6552.1763
$ws.Range("K136").Value = 10682.1819
$ws.Range("L136").Value = 19656.5289
$ws.Range("M136").Value = -8132.1819
$ws.Range("N136").Value = -24756.5289

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 251994.28
$ws.Range("I62").Value = 430146.25
$ws.Range("K62").Value = 430146.25
$ws.Range("M62").Value = -429522.25

# Row 65
$ws.Range("H65").Value = 251994.28
$ws.Range("I65").Value = 430146.25
$ws.Range("K65").Value = 2150731.25
$ws.Range("M65").Value = -2147611.25

# Row 113
$ws.Range("H113").Value = 1464.9395
$ws.Range("I113").Value = 849.6774
$ws.Range("J113").Value = 11001.5
$ws.Range("K113").Value = 2549.0322
$ws.Range("L113").Value = 33004.5
$ws.Range("M113").Value = -379.0322000000001
$ws.Range("N113").Value = -37344.5

# Row 132
$ws.Range("H132").Value = 11665.3545
$ws.Range("I132").Value = 13638.375
$ws.Range("K132").Value = 40915.125
$ws.Range("M132").Value = -38385.125
